$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CHECKLIST MATCHS")

# Clear the "vertical distributor check" and "no balls inside distributor" rows
$ws.Range("A3").Value = ""
$ws.Range("B3").Value = ""
$ws.Range("A4").Value = ""
$ws.Range("B4").Value = ""

# Remove "Récupérer la balise Catadioptre" line
$ws.Range("A26").Value = ""

# Rename "Récupérer le Robot" to "Récupérer le Robot 1" (two-phase release of the big robot)
$ws.Range("A27").Value = "Récupérer le Robot 1"

$null = $ws.Range("A19").Select()

$voyage = $wb.Worksheets.Item("CHECKLIST VOYAGE")
$null = $voyage.Range("D7").Select()

$testRobot = $wb.Worksheets.Item("CHECKLIST TEST ROBOT")
$null = $testRobot.Range("F11").Select()

$null = $ws.Activate()
